$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells (I0, IF) in row 1, matching the style of
# the existing header cells (bold, bordered, centered).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the data rows: column I is a constant 1 ("I0"), column J mirrors
# the existing "IP" value in column H ("IF").
for ($r = 2; $r -le 25; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
